# Apply updated timestamps / labels to the task order workbook.
$wb = $excel.ActiveWorkbook

# Rename worksheets (tab names) to new timestamp-based identifiers.
$wb.Worksheets.Item(1).Name = "GNG_TO-16512556372807684"
$wb.Worksheets.Item(2).Name = "NB_TO-16512556399811816"
$wb.Worksheets.Item(3).Name = "RS_TO-1651255639989003"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512556400449636"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512556401215255"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512556372434235.csv"
$ws1.Range("B3").Value = "GNG_stims-16512556372617638.csv"
$ws1.Range("B4").Value = "go_stims-1651255637263763.csv"
$ws1.Range("B5").Value = "GNG_stims-16512556372777646.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16512556388316116.csv"
$ws2.Range("B3").Value = "TB-16512556386112998.csv"
$ws2.Range("B4").Value = "OB-16512556382594776.csv"
$ws2.Range("B5").Value = "ZB-match_4-16512556376046185.csv"
$ws2.Range("B6").Value = "TB-16512556399678848.csv"
$ws2.Range("B7").Value = "ZB-match_9-16512556377324433.csv"
$ws2.Range("B8").Value = "OB-16512556379240463.csv"
$ws2.Range("B9").Value = "ZB-match_6-16512556373178973.csv"
$ws2.Range("B10").Value = "OB-1651255638014144.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512556400125527.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255639991002.csv"
$ws4.Range("B4").Value = "MM_stims-16512556400283859.csv"
$ws4.Range("B5").Value = "ZM_stims-16512556400135543.csv"
$ws4.Range("B6").Value = "MM_stims-1651255640043966.csv"
$ws4.Range("B7").Value = "ZM_stims-16512556400293305.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1651255640090309.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512556401056747.csv"
$ws5.Range("B4").Value = "SAT_stims-165125564007518.csv"
$ws5.Range("B5").Value = "SAT_stims-1651255640051075.csv"
